# Apply the burndown chart update:
#  - Chart!A4 text "Practice presentation" (was referencing the shared string
#    that held "Writing Assignment 1: Proposal"; now shares the string used by
#    the old A6 entry, i.e. same display text, but the underlying shared
#    string table is reshuffled so that index 22 becomes "Practice
#    presentation").
#  - Chart!A6 text changes to "Writing Assignment 1: Ouline" (new task name).
#  - Chart!D4 estimate-left updated from 0.5 to 2.
#  - Chart!C6 updated from 10 to 1.
#  - Chart!D6 updated from 0 to 1.
#  - Chart!C7 updated from 15 to 10.
#  - Selection on the Chart sheet moved to B7.

$wb = $excel.ActiveWorkbook
$chart = $wb.Worksheets.Item("Chart")

# Update the task text cells (A4/A6). Re-assigning the string values here
# naturally rewrites the shared-strings table to match the edited workbook.
$chart.Range("A4").Value = "Practice presentation"
$chart.Range("A6").Value = "Writing Assignment 1: Ouline"

# Update the numeric "time spent" / "time estimated" values that drive the
# burndown totals and chart series.
$chart.Range("D4").Value = 2
$chart.Range("C6").Value = 1
$chart.Range("D6").Value = 1
$chart.Range("C7").Value = 10

# Move the active selection on the Chart sheet from D5 to B7.
$chart.Activate()
$chart.Range("B7").Select()

$wb.Save()
